$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H6").Value = 487
$ws.Range("I6").Value = 149.66667
$ws.Range("K6").Value = 449.00001
$ws.Range("M6").Value = -337.00001
$ws.Range("H17").Value = 3249.5
$ws.Range("J17").Value = 3249.5
$ws.Range("L17").Value = 9748.5
$ws.Range("N17").Value = -10084.5
$ws.Range("H40").Value = 8399.799999999999
$ws.Range("I40").Value = 7999
$ws.Range("K40").Value = 7999
$ws.Range("M40").Value = -7824
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H64").Value = 7395.6
$ws.Range("I64").Value = 7491.5
$ws.Range("J64").Value = 7331.6665
$ws.Range("K64").Value = 7491.5
$ws.Range("L64").Value = 7331.6665
$ws.Range("M64").Value = -7243.5
$ws.Range("N64").Value = -7827.6665
$ws.Range("H67").Value = 7395.6
$ws.Range("I67").Value = 7491.5
$ws.Range("J67").Value = 7331.6665
$ws.Range("K67").Value = 7491.5
$ws.Range("L67").Value = 7331.6665
$ws.Range("M67").Value = -6633.5
$ws.Range("N67").Value = -9047.666499999999
$ws.Range("H74").Value = 10400
$ws.Range("J74").Value = 10400
$ws.Range("L74").Value = 10400
$ws.Range("N74").Value = -12272
$ws.Range("H77").Value = 10400
$ws.Range("J77").Value = 10400
$ws.Range("L77").Value = 52000
$ws.Range("N77").Value = -61360
$ws.Range("H138").Value = 3781.9062
$ws.Range("I138").Value = 2999.6
$ws.Range("J138").Value = 3926.7778
$ws.Range("K138").Value = 8998.799999999999
$ws.Range("L138").Value = 11780.3334
$ws.Range("M138").Value = -3858.799999999999
$ws.Range("N138").Value = -22060.3334

$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H35").Value = 5794.8
$ws.Range("I35").Value = 5794.8
$ws.Range("K35").Value = 5794.8
$ws.Range("M35").Value = -5388.8
$ws.Range("H61").Value = 1778.7
$ws.Range("I61").Value = 973.375
$ws.Range("K61").Value = 973.375
$ws.Range("M61").Value = -761.375
$ws.Range("H136").Value = 1778.7
$ws.Range("I136").Value = 973.375
$ws.Range("K136").Value = 2920.125
$ws.Range("M136").Value = -370.125

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H39").Value = 15021.5
$ws.Range("J39").Value = 15021.5
$ws.Range("L39").Value = 15021.5
$ws.Range("N39").Value = -15799.5
$ws.Range("H82").Value = 47512.2
$ws.Range("J82").Value = 74283
$ws.Range("L82").Value = 74283
$ws.Range("N82").Value = -75049
$ws.Range("H85").Value = 47512.2
$ws.Range("J85").Value = 74283
$ws.Range("L85").Value = 74283
$ws.Range("N85").Value = -76935
$ws.Range("H86").Value = 1874.5
$ws.Range("I86").Value = 1874.5
$ws.Range("K86").Value = 1874.5
$ws.Range("M86").Value = -751.5
$ws.Range("H89").Value = 1874.5
$ws.Range("I89").Value = 1874.5
$ws.Range("K89").Value = 9372.5
$ws.Range("M89").Value = -3756.5
$ws.Range("H94").Value = 641.6
$ws.Range("I94").Value = 641.6
$ws.Range("K94").Value = 641.6
$ws.Range("M94").Value = -190.6
$ws.Range("H134").Value = 6616.5713
$ws.Range("I134").Value = 1139.1428
$ws.Range("K134").Value = 3417.4284
$ws.Range("M134").Value = -882.4284000000002

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H7").Value = 675
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 1050
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 3150
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -3374
$ws.Range("H60").Value = 350
$ws.Range("I60").Value = 350
$ws.Range("K60").Value = 1050
$ws.Range("M60").Value = -799
$ws.Range("H63").Value = 3177.4
$ws.Range("I63").Value = 3177.4
$ws.Range("K63").Value = 9532.200000000001
$ws.Range("M63").Value = -8783.200000000001
$ws.Range("H64").Value = 2411.3333
$ws.Range("I64").Value = 2410
$ws.Range("J64").Value = 2412
$ws.Range("K64").Value = 7230
$ws.Range("L64").Value = 7236
$ws.Range("M64").Value = -6960
$ws.Range("N64").Value = -7776
$ws.Range("H66").Value = 3177.4
$ws.Range("I66").Value = 3177.4
$ws.Range("K66").Value = 28596.6
$ws.Range("M66").Value = -24852.6
$ws.Range("H67").Value = 2411.3333
$ws.Range("I67").Value = 2410
$ws.Range("J67").Value = 2412
$ws.Range("K67").Value = 7230
$ws.Range("L67").Value = 7236
$ws.Range("M67").Value = -6294
$ws.Range("N67").Value = -9108
$ws.Range("H102").Value = 7500
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H114").Value = 1533.6666
$ws.Range("I114").Value = 1
$ws.Range("J114").Value = 2300
$ws.Range("K114").Value = 3
$ws.Range("L114").Value = 6900
$ws.Range("M114").Value = 3251
$ws.Range("N114").Value = -13408
$ws.Range("H117").Value = 999
$ws.Range("I117").Value = 999
$ws.Range("K117").Value = 2997
$ws.Range("M117").Value = 445

$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0

$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 5155.8125
$ws.Range("I22").Value = 1599.2
$ws.Range("K22").Value = 1599.2
$ws.Range("M22").Value = -1304.2
$ws.Range("H27").Value = 5155.8125
$ws.Range("I27").Value = 1599.2
$ws.Range("K27").Value = 1599.2
$ws.Range("M27").Value = -1492.2
$ws.Range("H68").Value = 1900
$ws.Range("I68").Value = 1900
$ws.Range("K68").Value = 1900
$ws.Range("M68").Value = -1151
$ws.Range("H71").Value = 1900
$ws.Range("I71").Value = 1900
$ws.Range("K71").Value = 9500
$ws.Range("M71").Value = -5756
$ws.Range("H93").Value = 1974.625
$ws.Range("I93").Value = 2066.3333
$ws.Range("K93").Value = 2066.3333
$ws.Range("M93").Value = -818.3332999999998
